$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Duel Decks Merfolk vs. Goblins (DDT)'
$ws.Range("A2").Value = 'Akki Coalflinger'
$ws.Range("A3").Value = 'Allosaurus Rider'
$ws.Range("A4").Value = 'Ambush Commander'
$ws.Range("A5").Value = 'Boggart Shenanigans'
$ws.Range("A6").Value = 'Clickslither'
$ws.Range("A7").Value = 'Elvish Eulogist'
$ws.Range("A8").Value = 'Elvish Harbinger'
$ws.Range("A9").Value = 'Elvish Promenade'
$ws.Range("A10").Value = 'Elvish Warrior'
$ws.Range("A11").Value = 'Emberwilde Augur'
$ws.Range("A12").Value = 'Flamewave Invoker'
$ws.Range("A13").Value = 'Forest'
$ws.Range("A14").Value = 'Forest'
$ws.Range("A15").Value = 'Forest'
$ws.Range("A16").Value = 'Forest'
$ws.Range("A17").Value = 'Forgotten Cave'
$ws.Range("A18").Value = 'Gempalm Incinerator'
$ws.Range("A19").Value = 'Gempalm Strider'
$ws.Range("A20").Value = 'Giant Growth'
$ws.Range("A21").Value = 'Goblin Burrows'
$ws.Range("A22").Value = 'Goblin Cohort'
$ws.Range("A23").Value = 'Goblin Matron'
$ws.Range("A24").Value = 'Goblin Ringleader'
$ws.Range("A25").Value = 'Goblin Sledder'
$ws.Range("A26").Value = 'Goblin Warchief'
$ws.Range("A27").Value = 'Harmonize'
$ws.Range("A28").Value = 'Heedless One'
$ws.Range("A29").Value = 'Ib Halfheart, Goblin Tactician'
$ws.Range("A30").Value = 'Imperious Perfect'
$ws.Range("A31").Value = 'Llanowar Elves'
$ws.Range("A32").Value = 'Lys Alana Huntmaster'
$ws.Range("A33").Value = 'Mogg Fanatic'
$ws.Range("A34").Value = 'Mogg War Marshal'
$ws.Range("A35").Value = 'Moonglove Extract'
$ws.Range("A36").Value = 'Mountain'
$ws.Range("A37").Value = 'Mountain'
$ws.Range("A38").Value = 'Mountain'
$ws.Range("A39").Value = 'Mountain'
$ws.Range("A40").Value = 'Mudbutton Torchrunner'
$ws.Range("A41").Value = 'Raging Goblin'
$ws.Range("A42").Value = 'Reckless One'
$ws.Range("A43").Value = 'Siege-Gang Commander'
$ws.Range("A44").Value = 'Skirk Drill Sergeant'
$ws.Range("A45").Value = 'Skirk Fire Marshal'
$ws.Range("A46").Value = 'Skirk Prospector'
$ws.Range("A47").Value = 'Skirk Shaman'
$ws.Range("A48").Value = 'Slate of Ancestry'
$ws.Range("A49").Value = 'Spitting Earth'
$ws.Range("A50").Value = 'Stonewood Invoker'
$ws.Range("A51").Value = 'Sylvan Messenger'
$ws.Range("A52").Value = 'Tarfire'
$ws.Range("A53").Value = 'Tar Pitcher'
$ws.Range("A54").Value = 'Timberwatch Elf'
$ws.Range("A55").Value = 'Tranquil Thicket'
$ws.Range("A56").Value = 'Voice of the Woods'
$ws.Range("A57").Value = 'Wellwisher'
$ws.Range("A58").Value = 'Wildsize'
$ws.Range("A59").Value = 'Wirewood Herald'
$ws.Range("A60").Value = 'Wirewood Lodge'
$ws.Range("A61").Value = 'Wirewood Symbiote'
$ws.Range("A62").Value = 'Wood Elves'
$ws.Range("A63").Value = 'Wren''s Run Vanquisher'

$ws.Range("A64").ClearContents()
